# Reorder the worksheet tabs so "review_info" comes before "hotel_info",
# then add a new "State" column to the hotel_info sheet (inserted right
# after "Hotel_Name", before "City") with the value "Louisiana" for the
# single hotel row.

$wb = $excel.ActiveWorkbook

$hotelWs = $wb.Worksheets.Item("hotel_info")
$reviewWs = $wb.Worksheets.Item("review_info")

# Insert a new column C ("State") on hotel_info, shifting City (and
# everything after it) one column to the right.
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Range("C1").Value = "State"
$hotelWs.Range("C2").Value = "Louisiana"

# Move review_info so it precedes hotel_info in tab order. Done last so
# the worksheet references above aren't disturbed by the reorder.
$reviewWs.Move($hotelWs)
